$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated save_data: column G ("K") now holds strikeouts (K) recomputed
# from the game logs instead of the old Strike# count, so the per-appearance
# s_vals had to be recalculated and rewritten for every row of the sheet.
$kValues = @{
    2 = 2
    3 = 0
    4 = 2
    5 = 2
    7 = 0
    9 = 0
    10 = 0
    11 = 1
    13 = 0
    14 = 1
    15 = 2
    16 = 1
    17 = 1
    18 = 1
    19 = 3
    20 = 0
    21 = 0
    22 = 2
    23 = 1
    24 = 1
    25 = 1
    26 = 1
    27 = 0
    28 = 1
    29 = 1
    30 = 0
    31 = 2
    32 = 0
    33 = 0
    34 = 0
    35 = 1
    36 = 1
    37 = 2
    38 = 1
    39 = 0
    40 = 0
    41 = 3
    43 = 1
    44 = 1
    45 = 0
    46 = 1
    47 = 0
    48 = 0
    49 = 2
    50 = 2
    51 = 2
    52 = 3
    53 = 0
    54 = 0
    55 = 1
    56 = 1
    57 = 4
    58 = 0
    59 = 0
    63 = 2
}

foreach ($row in $kValues.Keys) {
    # Column 7 = G = "K"
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}

Write-Host "Updated $($kValues.Count) K values in column G"
